# Update 'F' column ('想去人数') values across all sheets as per upstream data refresh
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1416
$ws.Range("F5").Value = 6487
$ws.Range("F6").Value = 506
$ws.Range("F7").Value = 1054
$ws.Range("F8").Value = 21
$ws.Range("F9").Value = 3490
$ws.Range("F10").Value = 6724
$ws.Range("F12").Value = 1361
$ws.Range("F13").Value = 786
$ws.Range("F17").Value = 1123
$ws.Range("F19").Value = 117
$ws.Range("F21").Value = 181
$ws.Range("F23").Value = 1018
$ws.Range("F24").Value = 330
$ws.Range("F26").Value = 24
$ws.Range("F27").Value = 115
$ws.Range("F31").Value = 75
$ws.Range("F35").Value = 507
$ws.Range("F36").Value = 336
$ws.Range("F37").Value = 28
$ws.Range("F38").Value = 51
$ws.Range("F39").Value = 301
$ws.Range("F41").Value = 505
$ws.Range("F42").Value = 56
$ws.Range("F46").Value = 2

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F30").Value = 716
$ws.Range("F32").Value = 574
$ws.Range("F40").Value = 55

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F6").Value = 581
$ws.Range("F8").Value = 1171

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F6").Value = 581
$ws.Range("F10").Value = 6487
$ws.Range("F11").Value = 506
$ws.Range("F12").Value = 1054
$ws.Range("F13").Value = 21
$ws.Range("F15").Value = 6724
$ws.Range("F18").Value = 1361
$ws.Range("F24").Value = 1171
$ws.Range("F29").Value = 117
$ws.Range("F33").Value = 24
$ws.Range("F37").Value = 75
$ws.Range("F42").Value = 574
$ws.Range("F43").Value = 336
$ws.Range("F44").Value = 51
$ws.Range("F46").Value = 301
$ws.Range("F50").Value = 55
